# Applies the "updated required changes, pr lines, new liquidation features" edit:
#  - PRs sheet: new PR rows (p0123 / pr12323), drop the trailing
#    project_code/task_number/task_desc/wbl_percentage columns
#  - Payments sheet: add work_confirmation column + two payment rows
#  - DSA_Payments / Operational_Advances: unaffected data, headers stay the same text
#  - Timeline sheet: append new status-change history rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "PRs"
# ---------------------------------------------------------------------------
$prs = $wb.Worksheets.Item("PRs")

# Header row stays id..assigned_to (A..X); drop the old trailing
# project_code/task_number/task_desc/wbl_percentage columns (Y:AB).
$prs.Range("Y1:AB1").EntireColumn.Delete()

# Row 2
$prs.Cells.Item(2, 1).Value = 1
$prs.Cells.Item(2, 2).Value = "p0123"
$prs.Cells.Item(2, 3).Value = "2025-10-02"
$prs.Cells.Item(2, 4).Value = "admin"
$prs.Cells.Item(2, 5).Value = "CRLR"
$prs.Cells.Item(2, 6).Value = "Goods"
$prs.Cells.Item(2, 7).Value = "Rental Vehicle"
$prs.Cells.Item(2, 8).Value = "ssss"
$prs.Cells.Item(2, 9).Value = "Sedan Car"
$prs.Cells.Item(2, 10).Value = "fahad"
$prs.Cells.Item(2, 11).Value = "11111111111"
$prs.Cells.Item(2, 12).Value = "2025-10-01"
$prs.Cells.Item(2, 13).Value = "2025-10-02"
$prs.Cells.Item(2, 14).Value = 1
$prs.Cells.Item(2, 15).Value = "Islamabad"
$prs.Cells.Item(2, 16).Value = 2
$prs.Cells.Item(2, 17).Value = 0
$prs.Cells.Item(2, 18).Value = 0
$prs.Cells.Item(2, 19).Value = "Yes"
$prs.Cells.Item(2, 20).Value = 1
$prs.Cells.Item(2, 21).Value = ""
$prs.Cells.Item(2, 22).Value = "Submitted"
$prs.Cells.Item(2, 23).Value = "2025-10-02 15:54:37.236771"
$prs.Cells.Item(2, 24).Value = "admin"

# Row 3
$prs.Cells.Item(3, 1).Value = 2
$prs.Cells.Item(3, 2).Value = "pr12323"
$prs.Cells.Item(3, 3).Value = "2025-10-02"
$prs.Cells.Item(3, 4).Value = "admin"
$prs.Cells.Item(3, 5).Value = "CRLR"
$prs.Cells.Item(3, 6).Value = "Goods"
$prs.Cells.Item(3, 7).Value = "Rental Vehicle"
$prs.Cells.Item(3, 8).Value = ""
$prs.Cells.Item(3, 9).Value = "Sedan Car"
$prs.Cells.Item(3, 10).Value = "fahad"
$prs.Cells.Item(3, 11).Value = "1111111"
$prs.Cells.Item(3, 12).Value = "2025-10-01"
$prs.Cells.Item(3, 13).Value = "2025-10-02"
$prs.Cells.Item(3, 14).Value = 1
$prs.Cells.Item(3, 15).Value = "Islamabad"
$prs.Cells.Item(3, 16).Value = 1
$prs.Cells.Item(3, 17).Value = 100
$prs.Cells.Item(3, 18).Value = 100
$prs.Cells.Item(3, 19).Value = "Yes"
$prs.Cells.Item(3, 20).Value = 1
$prs.Cells.Item(3, 21).Value = ""
$prs.Cells.Item(3, 22).Value = "Completed"
$prs.Cells.Item(3, 23).Value = "2025-10-02 16:04:51.683825"
$prs.Cells.Item(3, 24).Value = "admin"

# ---------------------------------------------------------------------------
# Sheet "Payments"
# ---------------------------------------------------------------------------
$pay = $wb.Worksheets.Item("Payments")

# Header row: insert new "work_confirmation" column before work_order_yesno,
# and append the two new trailing columns (status, created_at).
$pay.Cells.Item(1, 1).Value = "id"
$pay.Cells.Item(1, 2).Value = "pr_id"
$pay.Cells.Item(1, 3).Value = "pr_number"
$pay.Cells.Item(1, 4).Value = "category"
$pay.Cells.Item(1, 5).Value = "po_number"
$pay.Cells.Item(1, 6).Value = "invoice_number"
$pay.Cells.Item(1, 7).Value = "wave_receipt"
$pay.Cells.Item(1, 8).Value = "work_confirmation"
$pay.Cells.Item(1, 9).Value = "work_order_yesno"
$pay.Cells.Item(1, 10).Value = "work_order_number"
$pay.Cells.Item(1, 11).Value = "actual_usd"
$pay.Cells.Item(1, 12).Value = "actual_pkr"
$pay.Cells.Item(1, 13).Value = "payment_date"
$pay.Cells.Item(1, 14).Value = "remarks"
$pay.Cells.Item(1, 15).Value = "status"
$pay.Cells.Item(1, 16).Value = "created_at"

# Row 2 (new payment tied to PR p0123)
$pay.Cells.Item(2, 1).Value = 1
$pay.Cells.Item(2, 2).Value = 1
$pay.Cells.Item(2, 3).Value = "p0123"
$pay.Cells.Item(2, 4).Value = "Rental Vehicle"
$pay.Cells.Item(2, 5).Value = "po123"
$pay.Cells.Item(2, 6).Value = "dd"
$pay.Cells.Item(2, 7).Value = "dd"
$pay.Cells.Item(2, 8).Value = "Yes"
$pay.Cells.Item(2, 11).Value = 0
$pay.Cells.Item(2, 12).Value = 0
$pay.Cells.Item(2, 13).Value = "2025-10-02"
$pay.Cells.Item(2, 15).Value = "In Process"
$pay.Cells.Item(2, 16).Value = "2025-10-02 10:57:17"

# Row 3 (new payment tied to PR pr12323)
$pay.Cells.Item(3, 1).Value = 2
$pay.Cells.Item(3, 2).Value = 2
$pay.Cells.Item(3, 3).Value = "pr12323"
$pay.Cells.Item(3, 4).Value = "Rental Vehicle"
$pay.Cells.Item(3, 5).Value = "p333"
$pay.Cells.Item(3, 6).Value = "dad"
$pay.Cells.Item(3, 7).Value = "dd"
$pay.Cells.Item(3, 8).Value = "Yes"
$pay.Cells.Item(3, 11).Value = 100
$pay.Cells.Item(3, 12).Value = 10
$pay.Cells.Item(3, 13).Value = "2025-10-02"
$pay.Cells.Item(3, 15).Value = "Completed"
$pay.Cells.Item(3, 16).Value = "2025-10-02 11:06:33"

# ---------------------------------------------------------------------------
# Sheet "Timeline": append the new status-change history rows
# ---------------------------------------------------------------------------
$tl = $wb.Worksheets.Item("Timeline")

$tl.Cells.Item(2, 1).Value = 1
$tl.Cells.Item(2, 2).Value = "PR"
$tl.Cells.Item(2, 3).Value = 1
$tl.Cells.Item(2, 5).Value = "Submitted"
$tl.Cells.Item(2, 6).Value = "admin"
$tl.Cells.Item(2, 7).Value = "2025-10-02 10:54:37"

$tl.Cells.Item(3, 1).Value = 2
$tl.Cells.Item(3, 2).Value = "Payment"
$tl.Cells.Item(3, 3).Value = 1
$tl.Cells.Item(3, 5).Value = "Pending"
$tl.Cells.Item(3, 6).Value = "admin"
$tl.Cells.Item(3, 7).Value = "2025-10-02 10:57:17"

$tl.Cells.Item(4, 1).Value = 3
$tl.Cells.Item(4, 2).Value = "Payment"
$tl.Cells.Item(4, 3).Value = 1
$tl.Cells.Item(4, 4).Value = "Pending"
$tl.Cells.Item(4, 5).Value = "In Process"
$tl.Cells.Item(4, 6).Value = "admin"
$tl.Cells.Item(4, 7).Value = "2025-10-02 10:57:47"

$tl.Cells.Item(5, 1).Value = 4
$tl.Cells.Item(5, 2).Value = "PR"
$tl.Cells.Item(5, 3).Value = 2
$tl.Cells.Item(5, 5).Value = "Submitted"
$tl.Cells.Item(5, 6).Value = "admin"
$tl.Cells.Item(5, 7).Value = "2025-10-02 11:04:51"

$tl.Cells.Item(6, 1).Value = 5
$tl.Cells.Item(6, 2).Value = "Payment"
$tl.Cells.Item(6, 3).Value = 2
$tl.Cells.Item(6, 5).Value = "Pending"
$tl.Cells.Item(6, 6).Value = "admin"
$tl.Cells.Item(6, 7).Value = "2025-10-02 11:06:33"

$tl.Cells.Item(7, 1).Value = 6
$tl.Cells.Item(7, 2).Value = "Payment"
$tl.Cells.Item(7, 3).Value = 2
$tl.Cells.Item(7, 4).Value = "Pending"
$tl.Cells.Item(7, 5).Value = "Completed"
$tl.Cells.Item(7, 6).Value = "admin"
$tl.Cells.Item(7, 7).Value = "2025-10-02 11:07:12"

$tl.Cells.Item(8, 1).Value = 7
$tl.Cells.Item(8, 2).Value = "PR"
$tl.Cells.Item(8, 3).Value = 2
$tl.Cells.Item(8, 4).Value = "Submitted"
$tl.Cells.Item(8, 5).Value = "Completed"
$tl.Cells.Item(8, 6).Value = "admin"
$tl.Cells.Item(8, 7).Value = "2025-10-02 11:07:12"
